# Normalize the "Recorded By" (column G) values so that "System" always
# appears first in the comma-separated list of recorders.
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com" -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
